$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the "Exported" timestamp text (A4)
# ------------------------------------------------------------------
$ws.Range("A4").Value = "Exported: 2018-04-20 20:19:47"

# ------------------------------------------------------------------
# 2. Preserve formatting for the trailing "EndOfData" marker row.
#    It currently lives at row 14; in the edited sheet it must move
#    down to row 16 (rows 12-13 become new data rows, pushing it down
#    by two).  Capture its format + row height before touching
#    anything else.
# ------------------------------------------------------------------
$endRowHeight = $ws.Rows.Item(14).RowHeight()
$ws.Range("A14:U14").Copy()
$ws.Range("A16:U16").PasteSpecial(-4122)
$ws.Range("A16").Value = "EndOfData"
# (Row height must be (re)applied AFTER the value is written, otherwise
# writing the cell value triggers an auto-fit that overrides it.)
$ws.Rows.Item(16).RowHeight = $endRowHeight

# Now clear out the old row 14 content/formatting completely so it
# collapses back to an untouched (default) row.
$ws.Range("A14:U14").Clear()
$ws.Rows.Item(14).AutoFit()

# ------------------------------------------------------------------
# 3. Build the two new data rows (12 = HELP 1, 13 = HELP 2) by
#    cloning the formatting of the existing row 11 ("Sample 002").
# ------------------------------------------------------------------
$ws.Range("A11:U11").Copy()
$ws.Range("A12:U12").PasteSpecial(-4122)
$ws.Range("A13:U13").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4. Row 9 ("Home"): flag it for listing, add keywords/description.
# ------------------------------------------------------------------
$ws.Range("M9").Value = 1
$ws.Range("P9").Value = "home,sample,pickles2"
$ws.Range("Q9").Value = "Home Page"

# ------------------------------------------------------------------
# 5. Row 10 ("Sample 001"): flag list + category-top.
# ------------------------------------------------------------------
$ws.Range("M10").Value = 1
$ws.Range("R10").Value = 1

# ------------------------------------------------------------------
# 6. Row 11 ("Sample 002"): flag list + category-top.
# ------------------------------------------------------------------
$ws.Range("M11").Value = 1
$ws.Range("R11").Value = 1

# ------------------------------------------------------------------
# 7. Row 12 ("HELP 1")
# ------------------------------------------------------------------
$ws.Range("C12").Value = "HELP 1"
$ws.Range("K12").Value = "/help/"
$ws.Range("M12").Value = 1

# ------------------------------------------------------------------
# 8. Row 13 ("HELP 2")
# ------------------------------------------------------------------
$ws.Range("C13").Value = "HELP 2"
$ws.Range("K13").Value = "/help/2.html"
$ws.Range("M13").Value = 1
